$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item($row, 3).Value = 'Bíobío'
$ws.Cells.Item($row, 4).Value = 44832
$ws.Cells.Item($row, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 'Fruta'
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = 'Otros'
$ws.Cells.Item($row, 9).Value = 100107002
$ws.Cells.Item($row, 10).Value = 'Chirimoya'
$ws.Cells.Item($row, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item($row, 12).Value = 'Primera'
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 25000
$ws.Cells.Item($row, 15).Value = 26000
$ws.Cells.Item($row, 16).Value = 25500
$ws.Cells.Item($row, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item($row, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item($row, 19).Value = 2550
$ws.Cells.Item($row, 20).Value = 10
